$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.009240865707397
$ws.Range("B1").Value = 2.119159460067749
$ws.Range("C1").Value = 6.367868900299072
$ws.Range("D1").Value = 1.567803382873535
$ws.Range("E1").Value = 1.362676978111267
